# "Generate Report for Handback"
# Updates the "Correspond Handback DateTime" column (H) on the ru-ru sheet
# for the rows whose handback report was (re)generated, stamping them with
# the latest handback run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ru-ru")

$newHandbackDateTime = "2016-03-24 08:48:51"
$rowsToStamp = @(9, 13, 14, 17, 18, 19, 23, 25, 31, 36, 40, 41, 43, 46, 47, 48, 52, 54, 59)

foreach ($row in $rowsToStamp) {
    $ws.Range("H$row").Value = $newHandbackDateTime
}
